# Lab 2 "AverageVelocities_60" edit: add a y60 standard-deviation column (C)
# next to the existing Position/y60 columns, and leave the selection where
# the author last clicked (D20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the existing "y60" column (header style + data
# style) onto the new column so the new cells pick up the same cell
# styles (border/alignment for the header, plain style for the data)
# instead of the worksheet's default style.
$ws.Range("B1:B13").Copy($ws.Range("C1:C13")) | Out-Null

# Header
$ws.Range("C1").Value = "y60StdDev"

# y60 standard deviation values, one per Position row (rows 2-13)
$stdDevValues = @(
    0.0583346997682705,
    0.054034435713544944,
    0.080567063690268947,
    0.038132402550844782,
    0.045587544586650437,
    0.030576111569808374,
    0.037411673467865142,
    0.090619556005781962,
    0.037447018645321642,
    0.075608984944315394,
    0.065207224420477064,
    0.034249418114256282
)

for ($i = 0; $i -lt $stdDevValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $stdDevValues[$i]
}

# The blank formatted rows below the table (15-16) pick up the wider
# column span once column C has data in the rows above, so touch and
# clear them to let that formatting catch up without leaving any content.
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 3).Value = $null
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 3).Value = $null

# Leave the cursor on D20, where the author's selection ended up.
$ws.Range("D20").Select() | Out-Null
